$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Export"
$ws.Range("D1").Value = "Status"

$ws.Range("C2").Value = $true
$ws.Range("D2").Value = "changes required"
